# Update the "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 612
    3  = 563
    5  = 28
    6  = 98
    7  = 56
    10 = 4913
    11 = 4618
    12 = 12
    16 = 170
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
